$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Title paragraph: "UI Comments " -> "UI " + "/" + "UX " (underlined)
#    + existing _GoBack bookmark + "Comments "
# ---------------------------------------------------------------------

# Remove the trailing "Comments " text (keep "UI " at the start); this
# also pulls the (end-of-paragraph) _GoBack bookmark left to sit right
# after "UI ".
$rComments = $d.Range(3, 12)
$rComments.Delete()

# Insert "/" then "UX " before the bookmark (i.e. right after "UI ").
$insSlash = $d.Range(3, 3)
$insSlash.InsertBefore("/")

$insUx = $d.Range(4, 4)
$insUx.InsertBefore("UX ")

# Underline just the "UX " run.
$rUx = $d.Range(4, 7)
$rUx.Font.Underline = 1

# Re-insert "Comments " after the bookmark.
$bm = $d.Bookmarks("_GoBack")
$insComments = $d.Range($bm.End, $bm.End)
$insComments.InsertAfter("Comments ")

# ---------------------------------------------------------------------
# 2. "The application visualization needs to be more clear"
#    -> "...more useable "
# ---------------------------------------------------------------------
$d.Content.Find.Execute("The application visualization needs to be more clear", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "The application visualization needs to be more useable ", 2) | Out-Null

# ---------------------------------------------------------------------
# 3. Remove the "application needs more ways of visualization..." bullet
# ---------------------------------------------------------------------
$pVis = $d.Paragraphs(5)
$pVis.Range.Delete()

# ---------------------------------------------------------------------
# 4. Remove the "design of the application needs to be more
#    inspiration" bullet
# ---------------------------------------------------------------------
$pDesign = $d.Paragraphs(5)
$pDesign.Range.Delete()

# ---------------------------------------------------------------------
# 5. "Tracking the financial status ..." -> "Adding a stock chart ..."
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Tracking the financial status for every month, and compare all months to each other, it will be great to appear on Stock chart ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Adding a stock chart to view the financial status over time-series ", 2) | Out-Null

# ---------------------------------------------------------------------
# 6. "Notify the user with his/her balance, ..." -> "Adding prediction
#    module of the future expenses for the user"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Notify the user with his/her balance, and the user can change the notification setting to be on a daily ,weekly  or monthly basis ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Adding prediction module of the future expenses for the user", 2) | Out-Null

# ---------------------------------------------------------------------
# 7. New bullet: "Adding audit log"
# ---------------------------------------------------------------------
$pAudit = $d.Paragraphs.Add()
$pAudit.Range.InsertAfter("Adding audit log")

# ---------------------------------------------------------------------
# 8. Two trailing empty "List Paragraph" paragraphs (no numbering)
# ---------------------------------------------------------------------
$pEmpty1 = $d.Paragraphs.Add()
$pEmpty1.Style = "List Paragraph"

$pEmpty2 = $d.Paragraphs.Add()
$pEmpty2.Style = "List Paragraph"

Write-Output "done"
